# Auto-generated edit script: updates leve-profit pricing/profit columns (H:N)
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect
# refreshed market-board prices pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1340.4546
$ws.Range("I2").Value = 1499.5
$ws.Range("J2").Value = 916.3333
$ws.Range("K2").Value = 1499.5
$ws.Range("L2").Value = 916.3333
$ws.Range("M2").Value = -1386.5
$ws.Range("N2").Value = -1142.3333
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H106").Value = 1992.5
$ws.Range("I106").Value = 1824.5834
$ws.Range("K106").Value = 1824.5834
$ws.Range("M106").Value = -1193.5834
$ws.Range("H137").Value = 18598.209
$ws.Range("I137").Value = 19842.834
$ws.Range("J137").Value = 14864.333
$ws.Range("K137").Value = 59528.50199999999
$ws.Range("L137").Value = 44592.999
$ws.Range("M137").Value = -56978.50199999999
$ws.Range("N137").Value = -49692.999
$ws.Range("H140").Value = 180000
$ws.Range("J140").Value = 180000
$ws.Range("L140").Value = 180000
$ws.Range("N140").Value = -190360
$ws.Range("H141").Value = 2141.8948
$ws.Range("I141").Value = 1984.0769
$ws.Range("J141").Value = 2483.8333
$ws.Range("K141").Value = 5952.2307
$ws.Range("L141").Value = 7451.499899999999
$ws.Range("M141").Value = -772.2307000000001
$ws.Range("N141").Value = -17811.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3634.5
$ws.Range("I2").Value = 2761.4
$ws.Range("K2").Value = 2761.4
$ws.Range("M2").Value = -2648.4
$ws.Range("H32").Value = 4020.45
$ws.Range("I32").Value = 3968.9824
$ws.Range("K32").Value = 3968.9824
$ws.Range("M32").Value = -3681.9824
$ws.Range("H61").Value = 3172.0715
$ws.Range("I61").Value = 2885.8
$ws.Range("J61").Value = 5557.6665
$ws.Range("K61").Value = 2885.8
$ws.Range("L61").Value = 5557.6665
$ws.Range("M61").Value = -2673.8
$ws.Range("N61").Value = -5981.6665
$ws.Range("H116").Value = 3634.5
$ws.Range("I116").Value = 2761.4
$ws.Range("K116").Value = 2761.4
$ws.Range("M116").Value = -467.4000000000001
$ws.Range("H132").Value = 32731.674
$ws.Range("I132").Value = 2357.561
$ws.Range("K132").Value = 7072.683000000001
$ws.Range("M132").Value = -4542.683000000001
$ws.Range("H136").Value = 3172.0715
$ws.Range("I136").Value = 2885.8
$ws.Range("J136").Value = 5557.6665
$ws.Range("K136").Value = 8657.400000000001
$ws.Range("L136").Value = 16672.9995
$ws.Range("M136").Value = -6107.400000000001
$ws.Range("N136").Value = -21772.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3634.5
$ws.Range("I3").Value = 2761.4
$ws.Range("K3").Value = 2761.4
$ws.Range("M3").Value = -2647.4
$ws.Range("H134").Value = 1300.6
$ws.Range("I134").Value = 1311.9231
$ws.Range("J134").Value = 1227
$ws.Range("K134").Value = 3935.7693
$ws.Range("L134").Value = 3681
$ws.Range("M134").Value = -1400.7693
$ws.Range("N134").Value = -8751

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4957.1934
$ws.Range("I31").Value = 3556.7058
$ws.Range("K31").Value = 3556.7058
$ws.Range("M31").Value = -3261.7058
$ws.Range("H34").Value = 4957.1934
$ws.Range("I34").Value = 3556.7058
$ws.Range("K34").Value = 3556.7058
$ws.Range("M34").Value = -3354.7058
$ws.Range("H74").Value = 1500000
$ws.Range("J74").Value = 1500000
$ws.Range("L74").Value = 1500000
$ws.Range("N74").Value = -1501748
$ws.Range("H77").Value = 1500000
$ws.Range("J77").Value = 1500000
$ws.Range("L77").Value = 4500000
$ws.Range("N77").Value = -4508736
$ws.Range("H88").Value = 31666.666
$ws.Range("J88").Value = 31666.666
$ws.Range("L88").Value = 31666.666
$ws.Range("N88").Value = -32478.666
$ws.Range("H91").Value = 31666.666
$ws.Range("J91").Value = 31666.666
$ws.Range("L91").Value = 31666.666
$ws.Range("N91").Value = -34474.666
$ws.Range("H99").Value = 6401.375
$ws.Range("I99").Value = 8500
$ws.Range("J99").Value = 2903.6667
$ws.Range("K99").Value = 8500
$ws.Range("L99").Value = 2903.6667
$ws.Range("M99").Value = -7002
$ws.Range("N99").Value = -5899.6667
$ws.Range("H126").Value = 6401.375
$ws.Range("I126").Value = 8500
$ws.Range("J126").Value = 2903.6667
$ws.Range("K126").Value = 25500
$ws.Range("L126").Value = 8711.000100000001
$ws.Range("M126").Value = -23030
$ws.Range("N126").Value = -13651.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 26916486
$ws.Range("I4").Value = 33438356
$ws.Range("K4").Value = 100315068
$ws.Range("M4").Value = -100314956
$ws.Range("H54").Value = 2499.5
$ws.Range("I54").Value = 999
$ws.Range("K54").Value = 2997
$ws.Range("M54").Value = -2438
$ws.Range("H55").Value = 5821669.5
$ws.Range("I55").Value = 1500175.6
$ws.Range("J55").Value = 9525807
$ws.Range("K55").Value = 4500526.800000001
$ws.Range("L55").Value = 28577421
$ws.Range("M55").Value = -4500349.800000001
$ws.Range("N55").Value = -28577775
$ws.Range("H59").Value = 1382.6666
$ws.Range("I59").Value = 566.3333
$ws.Range("J59").Value = 2199
$ws.Range("K59").Value = 1698.9999
$ws.Range("L59").Value = 6597
$ws.Range("M59").Value = -1158.9999
$ws.Range("N59").Value = -7677
$ws.Range("H97").Value = 1143.5
$ws.Range("I97").Value = 2077
$ws.Range("K97").Value = 6231
$ws.Range("M97").Value = -5735
$ws.Range("H129").Value = 1478.3793
$ws.Range("I129").Value = 1117.5555
$ws.Range("J129").Value = 1640.75
$ws.Range("K129").Value = 3352.6665
$ws.Range("L129").Value = 4922.25
$ws.Range("M129").Value = 1647.3335
$ws.Range("N129").Value = -14922.25
$ws.Range("H132").Value = 4153.6665
$ws.Range("I132").Value = 1949.5
$ws.Range("J132").Value = 4783.4287
$ws.Range("K132").Value = 17545.5
$ws.Range("L132").Value = 43050.85830000001
$ws.Range("M132").Value = -15015.5
$ws.Range("N132").Value = -48110.85830000001
$ws.Range("H138").Value = 7514
$ws.Range("I138").Value = 11995
$ws.Range("K138").Value = 35985
$ws.Range("M138").Value = -30845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 21923.75
$ws.Range("I24").Value = 19480
$ws.Range("J24").Value = 25996.666
$ws.Range("K24").Value = 19480
$ws.Range("L24").Value = 25996.666
$ws.Range("M24").Value = -19307
$ws.Range("N24").Value = -26342.666
$ws.Range("H102").Value = 45000
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4299.5713
$ws.Range("I122").Value = 3899.6
$ws.Range("K122").Value = 11698.8
$ws.Range("M122").Value = -9248.799999999999
$ws.Range("H136").Value = 2807.5652
$ws.Range("I136").Value = 2733.8235
$ws.Range("K136").Value = 8201.470499999999
$ws.Range("M136").Value = -5651.470499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 17247.25
$ws.Range("J37").Value = 17996.334
$ws.Range("L37").Value = 17996.334
$ws.Range("N37").Value = -18402.334
$ws.Range("H132").Value = 1403.35
$ws.Range("I132").Value = 1138.3334
$ws.Range("J132").Value = 2198.4
$ws.Range("K132").Value = 3415.0002
$ws.Range("L132").Value = 6595.200000000001
$ws.Range("M132").Value = -885.0001999999999
$ws.Range("N132").Value = -11655.2
$ws.Range("H136").Value = 6877.8184
$ws.Range("I136").Value = 5737.4287
$ws.Range("K136").Value = 17212.2861
$ws.Range("M136").Value = -14662.2861
